# Adds a new "2021" data column (R) to the 3.9.2 indicator sheet, mirroring
# the formatting of the existing last data column (Q) / the neighbouring
# data columns (N:P) for the body rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (year headers) ---------------------------------------------
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

# --- Row 5 --------------------------------------------------------------
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 1

# --- Row 6 ---------------------------------------------------------------
$ws.Range("P6").Copy($ws.Range("R6"))
$ws.Range("R6").Value = 2.2

# --- Row 7 ---------------------------------------------------------------
$ws.Range("P7").Copy($ws.Range("R7"))
$ws.Range("R7").Value = 1.7

# --- Row 8 (no data available -> "-") ------------------------------------
$ws.Range("P8").Copy($ws.Range("R8"))
$ws.Range("R8").Value = "-"

# --- Row 9 ---------------------------------------------------------------
$ws.Range("P9").Copy($ws.Range("R9"))
$ws.Range("R9").Value = 0.3

# --- Row 10 --------------------------------------------------------------
$ws.Range("P10").Copy($ws.Range("R10"))
$ws.Range("R10").Value = 1.1

# --- Row 11 (no data available -> "-") ------------------------------------
$ws.Range("P11").Copy($ws.Range("R11"))
$ws.Range("R11").Value = "-"

# --- Row 12 --------------------------------------------------------------
$ws.Range("P12").Copy($ws.Range("R12"))
$ws.Range("R12").Value = 0.9

# --- Row 13 --------------------------------------------------------------
$ws.Range("P13").Copy($ws.Range("R13"))
$ws.Range("R13").Value = 0.4

# --- Row 14 --------------------------------------------------------------
$ws.Range("Q14").Copy($ws.Range("R14"))
$ws.Range("R14").Value = 0.6

# Clear clipboard marquee / selection, move the active selection the same
# way the author's workbook ended up (one column further right than
# before, since a new rightmost data column was inserted).
$excel.CutCopyMode = $false
$ws.Range("S17").Select() | Out-Null
